$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The old row 33 (S33/T33 placeholder cells) is removed entirely; rows 34/35
# (the two footnote lines) shift up to become rows 33/34.
$ws.Rows.Item(33).Delete()

# New header row 37 (bold).
$ws.Range("A37").Value = "Filename"
$ws.Range("A37").Font.Bold = $true
$ws.Range("B37").Value = "Reads"
$ws.Range("B37").Font.Bold = $true
$ws.Range("C37").Value = "R1 mean Q30 to base"
$ws.Range("C37").Font.Bold = $true
$ws.Range("D37").Value = "R2 mean Q30 to base"
$ws.Range("D37").Font.Bold = $true
$ws.Range("E37").Value = "Sample Alias"
$ws.Range("E37").Font.Bold = $true
$ws.Range("F37").Value = "Sample Description"
$ws.Range("F37").Font.Bold = $true
$ws.Range("G37").Value = "Library Name"
$ws.Range("G37").Font.Bold = $true
$ws.Range("H37").Value = "Dilution Name"
$ws.Range("H37").Font.Bold = $true
$ws.Range("I37").Value = "Tag Barcode"
$ws.Range("I37").Font.Bold = $true
$ws.Range("J37").Value = "Mean insert size less adaptors"
$ws.Range("J37").Font.Bold = $true
$ws.Range("K37").Value = "Run Alias"
$ws.Range("K37").Font.Bold = $true
$ws.Range("L37").Value = "Lane"
$ws.Range("L37").Font.Bold = $true

# New data row 38.
$ws.Range("A38").Value = "1566_LIB18620_LDI16209_NoIndex_L002_R1_001.fastq.gz"
$ws.Range("A38").Font.Bold = $true
$ws.Range("B38").Value = 161482796
$ws.Range("B38").NumberFormat = "#,##0"
$ws.Range("C38").Value = 219
$ws.Range("D38").Value = 189
$ws.Range("E38").Value = "PRO973_S3_gDNA"
$ws.Range("F38").Value = "R.padi"
$ws.Range("G38").Value = "LIB18620"
$ws.Range("H38").Value = "LDI16209"
$ws.Range("I38").Value = "1: Index 2 (CGATGT)"
$ws.Range("J38").Value = 405
$ws.Range("K38").Value = "150825_SN790_0030_BH77HYBCXX"
$ws.Range("L38").Value = 2

# New highlighted section title on row 36 (bold, yellow fill) - added last,
# matching the authoring order (its shared string is appended after the
# table content).
$ws.Range("A36").Value = "R. padi GENOME read"
$ws.Range("A36").Font.Bold = $true
$ws.Range("A36").Interior.Color = 65535

$ws.Range("B30").Select()
